$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = 0.0026
$ws.Range("E4").Value = 118.12
$ws.Range("F4").Value = 0.00259
$ws.Range("H4").Value = 123.08

# Row 5
$ws.Range("C5").Value = 0.00336
$ws.Range("E5").Value = 117.01
$ws.Range("F5").Value = 0.00338
$ws.Range("H5").Value = 123.35

# Row 6
$ws.Range("C6").Value = 0.00263
$ws.Range("E6").Value = 115.15
$ws.Range("F6").Value = 0.00267
$ws.Range("H6").Value = 122.05

# Row 7
$ws.Range("C7").Value = 0.00339
$ws.Range("E7").Value = 115.1
$ws.Range("F7").Value = 0.00342
$ws.Range("H7").Value = 122.42

# Row 8
$ws.Range("E8").Value = 114.52
$ws.Range("F8").Value = 0.00239
$ws.Range("H8").Value = 121.97

# Row 9
$ws.Range("C9").Value = 0.00278
$ws.Range("E9").Value = 114.55
$ws.Range("H9").Value = 121.7

# Row 10
$ws.Range("C10").Value = 0.00265
$ws.Range("E10").Value = 115
$ws.Range("F10").Value = 0.00267
$ws.Range("H10").Value = 121.97

# Row 11
$ws.Range("C11").Value = 0.00311
$ws.Range("E11").Value = 115.68
$ws.Range("F11").Value = 0.00312
$ws.Range("H11").Value = 121.44

# Row 23
$ws.Range("C23").Value = 0.00268
$ws.Range("E23").Value = 115.64
$ws.Range("F23").Value = 0.00259
$ws.Range("H23").Value = 122.25

# Row 24
$ws.Range("C24").Value = 0.00341
$ws.Range("E24").Value = 115.64
$ws.Range("F24").Value = 0.00341
$ws.Range("H24").Value = 122.25

# Row 25
$ws.Range("C25").Value = 0.0026
$ws.Range("E25").Value = 115.64
$ws.Range("F25").Value = 0.00267
$ws.Range("H25").Value = 122.25

# Row 26
$ws.Range("C26").Value = 0.00336
$ws.Range("E26").Value = 115.64
$ws.Range("F26").Value = 0.00344
$ws.Range("H26").Value = 122.25

# Row 27
$ws.Range("C27").Value = 0.00231
$ws.Range("E27").Value = 115.64
$ws.Range("F27").Value = 0.00239
$ws.Range("H27").Value = 122.25

# Row 28
$ws.Range("C28").Value = 0.00274
$ws.Range("E28").Value = 115.64
$ws.Range("H28").Value = 122.25

# Row 29
$ws.Range("C29").Value = 0.00264
$ws.Range("E29").Value = 115.64
$ws.Range("F29").Value = 0.00266
$ws.Range("H29").Value = 122.25

# Row 30
$ws.Range("C30").Value = 0.00312
$ws.Range("E30").Value = 115.64
$ws.Range("F30").Value = 0.00309
$ws.Range("H30").Value = 122.25
